$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in Departamento name
$ws.Range("A7").Value = "Arequipa"

# Convert Albergue column (E) from numeric 1/0 flags to "Si"/"No" text
$ws.Range("E2").Value = "Si"
$ws.Range("E3").Value = "No"
$ws.Range("E4").Value = "Si"
$ws.Range("E5").Value = "No"
$ws.Range("E6").Value = "Si"
$ws.Range("E7").Value = "Si"
$ws.Range("E8").Value = "No"
$ws.Range("E9").Value = "No"
